$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume text values are preserved as text (matches source formatting)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "97.883.82"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "3.411.92"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "255.67"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").Value = "656.33"
$ws.Range("E6").Value = "  +4.81%  "
$ws.Range("D7").Value = "1.48"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("D8").Value = "0.437"
$ws.Range("E8").Value = "  +4.94%  "
$ws.Range("D9").Value = "1.08"
$ws.Range("E9").Value = "  +4.22%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "3.409.48"
$ws.Range("E11").Value = "  +2.96%  "
$ws.Range("E12").Value = "  +4.34%  "
$ws.Range("D13").Value = "41.97"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").Value = "6.43"
$ws.Range("E14").Value = "  +19.12%  "
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").Value = "97.658.61"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").Value = "4.055.15"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "8.63"
$ws.Range("E18").Value = "  +34.50%  "
$ws.Range("D19").Value = "3.414.21"
$ws.Range("E19").Value = "  +2.96%  "
$ws.Range("D20").Value = "17.67"
$ws.Range("E20").Value = "  +12.04%  "
$ws.Range("D21").Value = "0.495"
$ws.Range("E21").Value = "  +43.81%  "
$ws.Range("D22").Value = "3.47"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "10.76"
$ws.Range("E23").Value = "  +13.54%  "
$ws.Range("D24").Value = "508.09"
$ws.Range("E24").Value = "  +4.05%  "
$ws.Range("D25").Value = "0.0000207"
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("D26").Value = "6.23"
$ws.Range("E26").Value = "  +6.01%  "
$ws.Range("D27").Value = "99.38"
$ws.Range("E27").Value = "  +11.21%  "
$ws.Range("D28").Value = "12.91"
$ws.Range("E28").Value = "  +6.00%  "
$ws.Range("D29").Value = "3.595.87"
$ws.Range("E29").Value = "  +3.14%  "
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("E31").Value = "  +5.22%  "
$ws.Range("D32").Value = "11.47"
$ws.Range("E32").Value = "  +7.86%  "
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D35").Value = "0.574"
$ws.Range("E35").Value = "  +18.02%  "
$ws.Range("D36").Value = "29.77"
$ws.Range("E36").Value = "  +5.85%  "
$ws.Range("E37").Value = "  +16.70%  "
$ws.Range("D38").Value = "7.86"
$ws.Range("E38").Value = "  +6.82%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "533.41"
$ws.Range("E39").Value = "  +7.26%  "
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "1.44"
$ws.Range("E40").Value = "  +14.41%  "
$ws.Range("D41").Value = "0.154"
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("D42").Value = "24.74"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "0.868"
$ws.Range("E43").Value = "  +9.79%  "
$ws.Range("E44").Value = "  -4.57%  "
$ws.Range("D45").Value = "0.0421"
$ws.Range("E45").Value = "  +22.13%  "
$ws.Range("D46").Value = "3.33"
$ws.Range("E46").Value = "  +5.06%  "
$ws.Range("D47").Value = "5.52"
$ws.Range("E47").Value = "  +15.32%  "
$ws.Range("D48").Value = "8.32"
$ws.Range("E48").Value = "  +13.32%  "
$ws.Range("D49").Value = "1.61"
$ws.Range("E49").Value = "  +13.70%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("E51").Value = "  +5.99%  "
